$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item(1)   # TEST_CASES
$wsSteps     = $wb.Worksheets.Item(2)   # STEPS
$wsParams    = $wb.Worksheets.Item(3)   # PARAMETERS
$wsDatasets  = $wb.Worksheets.Item(4)   # DATASETS

# ---------------------------------------------------------------
# 1) TEST_CASES sheet: remove the trailing blank rows 12:16
# ---------------------------------------------------------------
$wsTestCases.Rows("12:16").Delete()

# ---------------------------------------------------------------
# 2) TEST_CASES sheet: fill in TC_PATH (col D) / TC_NUM (col E)
#    for the remaining data rows (2:11).
#    The write order below reproduces the original shared-string
#    insertion order (row 4 is populated last).
# ---------------------------------------------------------------
$wsTestCases.Range("D2").Value  = "path/row1"
$wsTestCases.Range("D3").Value  = "path/row2"
$wsTestCases.Range("D5").Value  = "path/row4"
$wsTestCases.Range("D6").Value  = "path/row5"
$wsTestCases.Range("D7").Value  = "path/row6"
$wsTestCases.Range("D8").Value  = "path/row7"
$wsTestCases.Range("D9").Value  = "path/row8"
$wsTestCases.Range("D10").Value = "path/row9"
$wsTestCases.Range("D11").Value = "path/row10"
$wsTestCases.Range("D4").Value  = "path/row3"

$wsTestCases.Range("E2").Value  = 11
$wsTestCases.Range("E3").Value  = 12
$wsTestCases.Range("E4").Value  = 13
$wsTestCases.Range("E5").Value  = 14
$wsTestCases.Range("E6").Value  = 15
$wsTestCases.Range("E7").Value  = 16
$wsTestCases.Range("E8").Value  = 17
$wsTestCases.Range("E9").Value  = 18
$wsTestCases.Range("E10").Value = 19
$wsTestCases.Range("E11").Value = 20

# Rows 3:11 need their D/E formatting aligned on the same
# border style already used by row 2 (D2:E2). Re-use that
# formatting via a format-only copy/paste so the existing
# style is reused instead of creating new ones.
$wsTestCases.Range("D2:E2").Copy()
$wsTestCases.Range("D3:E3").PasteSpecial(-4122)
$wsTestCases.Range("D4:E4").PasteSpecial(-4122)
$wsTestCases.Range("D5:E5").PasteSpecial(-4122)
$wsTestCases.Range("D6:E6").PasteSpecial(-4122)
$wsTestCases.Range("D7:E7").PasteSpecial(-4122)
$wsTestCases.Range("D8:E8").PasteSpecial(-4122)
$wsTestCases.Range("D9:E9").PasteSpecial(-4122)
$wsTestCases.Range("D10:E10").PasteSpecial(-4122)
$wsTestCases.Range("D11:E11").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 3) TEST_CASES sheet: widen column D slightly
# ---------------------------------------------------------------
$wsTestCases.Columns.Item(4).ColumnWidth = 10.3

# ---------------------------------------------------------------
# 4) Update the remembered selections on STEPS and DATASETS,
#    without leaving them as the active sheet/tab.
# ---------------------------------------------------------------
$wsSteps.Range("B34").Select()
$wsDatasets.Range("H7").Select()

# ---------------------------------------------------------------
# 5) Finally activate TEST_CASES with its new selection - this
#    becomes the active sheet/tab saved with the workbook.
# ---------------------------------------------------------------
$wsTestCases.Range("C18").Select()
